$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 5 corresponds to the "age_groups" variable.
# Update the description from "Five age groups" to "Four age groups"
$ws.Range("B5").Value = "Four age groups"

# Update the Levels column with the new four-group age breakdown
$ws.Range("D5").Value = "[60, 65], (65, 70], (70, 75],  (75, 80]"
